# Restructure the ticket importer (Data sheet, columns I & J) so the
# layout makes more sense:
#   - Column I header becomes "Open" (was "Additional Comments"); it keeps
#     the optional/cyan header style and gets a new comment explaining the
#     1/0 open flag.
#   - Column J header becomes "Additional Comments" (was "Open"); it takes
#     on the optional/cyan header style (matching column I) and inherits
#     the "you can add more comment columns" explanation that used to live
#     on I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$colI = $ws.Range("I1")
$colJ = $ws.Range("J1")

# --- Swap the header text -------------------------------------------------
$i1Value = $colI.Value2
$j1Value = $colJ.Value2
$colI.Value = $j1Value
$colJ.Value = $i1Value

# --- Give J1 the same (optional-column / cyan) formatting as I1 ----------
$colI.Copy()
$colJ.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the explanatory comments on the two header cells -------------
# The "additional comments" explanation moves from I1 down to J1 (its new
# home), and I1 gets a fresh explanation of the 1/0 "Open" flag.
$additionalCommentsText = $colI.Comment.Text()
$null = $colJ.Comment.Text($additionalCommentsText)
$null = $colI.Comment.Text("1 for open, 0 for closed. Defaults to open if left blank.")
